$d = $word.ActiveDocument
$d.Content.Find.Execute("Kim Phung Tran", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Phung Tran", 2)
